$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metro_budget")
$ws.Range("C98").Style = "Percent"
$ws.Range("C98").NumberFormat = "0.00%"
Write-Host "done"
